$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 208
$ws.Range("I2").Value = 583
$ws.Range("J2").Value = 2380
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 679
$ws.Range("M2").Value = 41
$ws.Range("N2").Value = 409
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 4
$ws.Range("R2").Value = 32
$ws.Range("S2").Value = 242
$ws.Range("T2").Value = 420
$ws.Range("U2").Value = 31
$ws.Range("V2").Value = 3883
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3765
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 54
$ws.Range("AA2").Value = 27
